# feat: add 2022-Q1 data
#
# The previous "总计" (totals) sheet becomes the new "2022-Q1" sheet (it
# keeps the same underlying sheet/rId, just like Excel does when you
# rename a tab), and a brand-new "总计" sheet is appended after it with
# the refreshed roll-up numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the old "总计" sheet into "2022-Q1" and add a fresh "总计"
#    sheet right after it, so tab order stays 2020-Q4, 2021-Q1,
#    2022-Q1, 总计.
# ---------------------------------------------------------------------
$quarterSheet = $wb.Worksheets.Item("总计")
$quarterSheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add($null, $quarterSheet)
$totalSheet.Name = "总计"

# Match the page-setup margins used throughout the rest of the workbook
# (a brand-new sheet otherwise gets Excel's generic defaults).
$totalSheet.PageSetup.LeftMargin = 72 * 0.75
$totalSheet.PageSetup.RightMargin = 72 * 0.75
$totalSheet.PageSetup.TopMargin = 72 * 1
$totalSheet.PageSetup.BottomMargin = 72 * 1
$totalSheet.PageSetup.HeaderMargin = 72 * 0.5
$totalSheet.PageSetup.FooterMargin = 72 * 0.5
$totalSheet.Outline.SummaryRow = 1
$totalSheet.Outline.SummaryColumn = 1

# ---------------------------------------------------------------------
# 2. Build the "2022-Q1" holdings sheet. Start from the "2021-Q1" sheet
#    layout (same header/format/border styling) then overwrite with the
#    2022-Q1 numbers.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q1")
$template.Range("B1:H1").Copy($quarterSheet.Range("B1:H1"))
$template.Range("A2:H3").Copy($quarterSheet.Range("A2:H3"))

# Headers (D1 differs from the template: "基金规模" instead of "基金金额")
$quarterSheet.Range("B1").Value = "基金代码"
$quarterSheet.Range("C1").Value = "基金名称"
$quarterSheet.Range("D1").Value = "基金规模"
$quarterSheet.Range("E1").Value = "股票总仓位"
$quarterSheet.Range("F1").Value = "仓位占比"
$quarterSheet.Range("G1").Value = "持有市值(亿元)"
$quarterSheet.Range("H1").Value = "仓位排名"

# Row 2 - numeric-looking values (fund code / amounts) must stay TEXT,
# so they are entered with a leading apostrophe (forces text, like
# typing into Excel by hand) and the style is reset back to Normal
# afterwards so no stray "number stored as text" formatting lingers.
$quarterSheet.Range("B2").Value = "'910021"
$quarterSheet.Range("B2").Style = "Normal"
$quarterSheet.Range("C2").Value = "东方红启华三年持有期混合型证券投资基金A"
$quarterSheet.Range("D2").Value = "'5.93"
$quarterSheet.Range("D2").Style = "Normal"
$quarterSheet.Range("E2").Value = "'86.11"
$quarterSheet.Range("E2").Style = "Normal"
$quarterSheet.Range("F2").Value = "'3.14"
$quarterSheet.Range("F2").Style = "Normal"
$quarterSheet.Range("G2").Value = "'0.1862"
$quarterSheet.Range("G2").Style = "Normal"
$quarterSheet.Range("H2").Value = 8

# Row 3
$quarterSheet.Range("B3").Value = "'011313"
$quarterSheet.Range("B3").Style = "Normal"
$quarterSheet.Range("C3").Value = "东方红启华三年持有期混合型证券投资基金B"
$quarterSheet.Range("D3").Value = "'0.97"
$quarterSheet.Range("D3").Style = "Normal"
$quarterSheet.Range("E3").Value = "'86.11"
$quarterSheet.Range("E3").Style = "Normal"
$quarterSheet.Range("F3").Value = "'3.14"
$quarterSheet.Range("F3").Style = "Normal"
$quarterSheet.Range("G3").Value = "'0.0305"
$quarterSheet.Range("G3").Style = "Normal"
$quarterSheet.Range("H3").Value = 8

# ---------------------------------------------------------------------
# 3. Build the new "总计" roll-up sheet: same header as before, plus a
#    new 2022-Q1 row on top, with the older rows shifted down.
# ---------------------------------------------------------------------
$template.Range("B1:D1").Copy($totalSheet.Range("B1:D1"))
$template.Range("A2:D2").Copy($totalSheet.Range("A2:D4"))

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.22

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q1"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 2.64

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2020-Q4"
$totalSheet.Range("C4").Value = 5
$totalSheet.Range("D4").Value = 0.41

# Restore the original active tab (the first sheet), mirroring the
# source workbook instead of leaving focus on the brand-new "总计" tab.
$wb.Worksheets.Item("2020-Q4").Activate()
